$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")

# Update HLSR data table with 2021 figures (rows shift down, 2015 dropped)
$ws1.Cells.Item(2,1).Value = 2021
$ws1.Cells.Item(2,2).Value = 662.7570079770649
$ws1.Cells.Item(2,3).Value = 7941263394.782457
$ws1.Cells.Item(2,4).Value = 11982164.351640126
$ws1.Cells.Item(2,5).Value = 0.6015846905886887
$ws1.Cells.Item(2,6).Value = 125.39389272546444
$ws1.Cells.Item(2,7).Value = 454.3177064337384
$ws1.Cells.Item(2,8).Value = -0.2525482205560857
$ws1.Cells.Item(2,9).Value = -0.04846991115653043
$ws1.Cells.Item(2,10).Value = 0.2730320737899432
$ws1.Cells.Item(2,11).Value = 0.2531685915640278
$ws1.Cells.Item(2,12).Value = -0.08379240839397273
$ws1.Cells.Item(2,13).Value = -0.24480503428890987
$ws1.Cells.Item(2,14).Value = 94.78111966362835
$ws1.Cells.Item(2,15).Value = 61.747354520507066
$ws1.Cells.Item(3,1).Value = 2020
$ws1.Cells.Item(3,2).Value = 886.6886482899803
$ws1.Cells.Item(3,3).Value = 8345782742.860617
$ws1.Cells.Item(3,4).Value = 9412303.584754176
$ws1.Cells.Item(3,5).Value = 0.48005088432504983
$ws1.Cells.Item(3,6).Value = 136.86187920104496
$ws1.Cells.Item(3,7).Value = 601.58995631804
$ws1.Cells.Item(3,8).Value = 1.2155110234593454
$ws1.Cells.Item(3,9).Value = -0.04100983702295247
$ws1.Cells.Item(3,10).Value = -0.5671471941134104
$ws1.Cells.Item(3,11).Value = -0.507170798339046
$ws1.Cells.Item(3,12).Value = 0.08629515675576349
$ws1.Cells.Item(3,13).Value = 1.2209110515285637
$ws1.Cells.Item(3,14).Value = 99.60916714554911
$ws1.Cells.Item(3,15).Value = 48.50416245733625
$ws1.Cells.Item(4,1).Value = 2019
$ws1.Cells.Item(4,2).Value = 400.21856759055345
$ws1.Cells.Item(4,3).Value = 8702678155.688616
$ws1.Cells.Item(4,4).Value = 21744813.610426877
$ws1.Cells.Item(4,5).Value = 0.9740715093731497
$ws1.Cells.Item(4,6).Value = 125.9895879585664
$ws1.Cells.Item(4,7).Value = 270.8753040352471
$ws1.Cells.Item(4,8).Value = -0.002474912027855214
$ws1.Cells.Item(4,9).Value = 0.014080730999245228
$ws1.Cells.Item(4,10).Value = 0.016596718445203384
$ws1.Cells.Item(4,11).Value = 0.010175199012249392
$ws1.Cells.Item(4,12).Value = 0.0063459591994667885
$ws1.Cells.Item(4,13).Value = -0.0018454097782281131
$ws1.Cells.Item(4,14).Value = 103.86881012035278
$ws1.Cells.Item(4,15).Value = 112.0569436023124
$ws1.Cells.Item(5,1).Value = 2018
$ws1.Cells.Item(5,2).Value = 401.2115308339285
$ws1.Cells.Item(5,3).Value = 8581839581.068909
$ws1.Cells.Item(5,4).Value = 21389812.90799727
$ws1.Cells.Item(5,5).Value = 0.9642599722558998
$ws1.Cells.Item(5,6).Value = 125.19510493070321
$ws1.Cells.Item(5,7).Value = 271.3761041514256
$ws1.Cells.Item(5,8).Value = -0.03610063789536211
$ws1.Cells.Item(5,9).Value = 0.015247568031038217
$ws1.Cells.Item(5,10).Value = 0.05327133510523696
$ws1.Cells.Item(5,11).Value = 0.04917543795846724
$ws1.Cells.Item(5,12).Value = -0.0004491920743003863
$ws1.Cells.Item(5,13).Value = -0.030649484459494225
$ws1.Cells.Item(5,14).Value = 102.4265691529347
$ws1.Cells.Item(5,15).Value = 110.22752834938696
$ws1.Cells.Item(6,1).Value = 2017
$ws1.Cells.Item(6,2).Value = 416.2379877064119
$ws1.Cells.Item(6,3).Value = 8452952611.0685
$ws1.Cells.Item(6,4).Value = 20307979.715274524
$ws1.Cells.Item(6,5).Value = 0.9190645695368166
$ws1.Cells.Item(6,6).Value = 125.25136685198841
$ws1.Cells.Item(6,7).Value = 279.95663054876223
$ws1.Cells.Item(6,8).Value = -0.035969301290042055
$ws1.Cells.Item(6,9).Value = 0.00888268416717164
$ws1.Cells.Item(6,10).Value = 0.04652547425847908
$ws1.Cells.Item(6,11).Value = 0.0421637732847191
$ws1.Cells.Item(6,12).Value = 0.011290411509540688
$ws1.Cells.Item(6,13).Value = -0.03902808574171002
$ws1.Cells.Item(6,14).Value = 100.88826841671717
$ws1.Cells.Item(6,15).Value = 104.65254742584791
$ws1.Cells.Item(7,1).Value = 2016
$ws1.Cells.Item(7,2).Value = 431.7683952009114
$ws1.Cells.Item(7,3).Value = 8378528786.076228
$ws1.Cells.Item(7,4).Value = 19405146.090365212
$ws1.Cells.Item(7,5).Value = 0.8818811333655217
$ws1.Cells.Item(7,6).Value = 123.85301534207888
$ws1.Cells.Item(7,7).Value = 291.32654804468666
$ws1.Cells.Item(7,8).Value = -0.019204070002918572
$ws1.Cells.Item(7,9).Value = 0.006568316112772976
$ws1.Cells.Item(7,10).Value = 0.026277011687607876
$ws1.Cells.Item(7,11).Value = 0.018155459390721607
$ws1.Cells.Item(7,12).Value = 0.016252900049103847
$ws1.Cells.Item(7,13).Value = -0.02734775827991065
$ws1.Cells.Item(7,14).Value = 100
$ws1.Cells.Item(7,15).Value = 100
# Make the data sheet the active/selected tab with new selection range
$ws1.Activate()
$ws1.Range("A2:O7").Select()
